{"js": "// Insert three new paragraphs (\"Modification 1\", \"Modification 2\",\n// \"Modification 3\") right after the second paragraph in the body (the\n// empty paragraph that precedes the paragraph holding the \"_GoBack\"\n// bookmark), matching the target diff.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// The second paragraph (index 1) is the empty paragraph just before the\n// trailing bookmark paragraph - the new paragraphs are inserted after it.\nconst anchorParagraph = paragraphs.items[1];\n\nlet lastParagraph = anchorParagraph;\nconst newTexts = [\"Modification 1\", \"Modification 2\", \"Modification 3\"];\nfor (const text of newTexts) {\n  lastParagraph = lastParagraph.insertParagraph(text, \"After\");\n}\n\nawait context.sync();\n", "ps1": "# Insert three new paragraphs (\"Modification 1\", \"Modification 2\",\n# \"Modification 3\") right after the second paragraph in the document body\n# (the empty paragraph that precedes the paragraph holding the \"_GoBack\"\n# bookmark), matching the target diff.\n\n$d = $word.ActiveDocument\n\n# The second paragraph is the empty paragraph just before the trailing\n# bookmark paragraph - the new paragraphs are inserted after it.\n$anchorIndex = 2\n\n$texts = @(\"Modification 1\", \"Modification 2\", \"Modification 3\")\nforeach ($text in $texts) {\n    $anchorParagraph = $d.Paragraphs.Item($anchorIndex)\n    $anchorParagraph.Range.InsertParagraphAfter()\n    $anchorIndex = $anchorIndex + 1\n    $d.Paragraphs.Item($anchorIndex).Range.InsertAfter($text)\n}\n"}
